$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update comparison values for the three existing model rows (3-5)
$ws.Range("B3").Value = 75
$ws.Range("C3").Value = 66.7

$ws.Range("B4").Value = 93.8
$ws.Range("C4").Value = 13.3

$ws.Range("B5").Value = 94
$ws.Range("C5").Value = 54.7

# Insert a new row before the old "CXR + CAD" row (row 6), shifting the
# remaining rows down, then populate it with the new model data.
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = "TPP Optimized with Urine LAM (Parallel)"
$ws.Range("B6").Value = 85.8
$ws.Range("C6").Value = 66
